# Auto-generated edit script applying numeric cell updates per the commit diff.
# Workbook contains 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), each backed
# by an Excel Table (ListObject) spanning A1:N141 with crafting-leve profit data.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 592.1754
$ws.Range("I15").Value = 592.1754
$ws.Range("K15").Value = 1776.5262
$ws.Range("M15").Value = -1607.5262
$ws.Range("H33").Value = 325.8889
$ws.Range("I33").Value = 130
$ws.Range("K33").Value = 130
$ws.Range("M33").Value = 99
$ws.Range("H40").Value = 1924.875
$ws.Range("I40").Value = 1914.1428
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1914.1428
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1739.1428
$ws.Range("N40").Value = -2350
$ws.Range("H62").Value = 6245.154
$ws.Range("I62").Value = 4359.5
$ws.Range("J62").Value = 7083.222
$ws.Range("K62").Value = 4359.5
$ws.Range("L62").Value = 7083.222
$ws.Range("M62").Value = -3735.5
$ws.Range("N62").Value = -8331.222
$ws.Range("H65").Value = 6245.154
$ws.Range("I65").Value = 4359.5
$ws.Range("J65").Value = 7083.222
$ws.Range("K65").Value = 21797.5
$ws.Range("L65").Value = 35416.11
$ws.Range("M65").Value = -18677.5
$ws.Range("N65").Value = -41656.11
$ws.Range("H113").Value = 5173
$ws.Range("I113").Value = 4712
$ws.Range("K113").Value = 4712
$ws.Range("M113").Value = -1458
$ws.Range("H116").Value = 4529.625
$ws.Range("I116").Value = 3069.6
$ws.Range("K116").Value = 3069.6
$ws.Range("M116").Value = 372.4000000000001
$ws.Range("H137").Value = 8052.8096
$ws.Range("I137").Value = 1613.909
$ws.Range("K137").Value = 4841.727000000001
$ws.Range("M137").Value = -2291.727000000001
$ws.Range("H138").Value = 5106.843
$ws.Range("I138").Value = 5972.2085
$ws.Range("J138").Value = 4655.3477
$ws.Range("K138").Value = 17916.6255
$ws.Range("L138").Value = 13966.0431
$ws.Range("M138").Value = -12776.6255
$ws.Range("N138").Value = -24246.0431

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1817.4166
$ws.Range("I2").Value = 2058
$ws.Range("J2").Value = 1095.6666
$ws.Range("K2").Value = 2058
$ws.Range("L2").Value = 1095.6666
$ws.Range("M2").Value = -1945
$ws.Range("N2").Value = -1321.6666
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 150
$ws.Range("K5").Value = 150
$ws.Range("M5").Value = -38
$ws.Range("H61").Value = 2399.6667
$ws.Range("I61").Value = 2399
$ws.Range("K61").Value = 2399
$ws.Range("M61").Value = -2187
$ws.Range("H95").Value = 52999.332
$ws.Range("J95").Value = 52999.332
$ws.Range("L95").Value = 52999.332
$ws.Range("N95").Value = -58491.332
$ws.Range("H97").Value = 587.5263
$ws.Range("I97").Value = 575.5333000000001
$ws.Range("K97").Value = 575.5333000000001
$ws.Range("M97").Value = -79.53330000000005
$ws.Range("H116").Value = 1817.4166
$ws.Range("I116").Value = 2058
$ws.Range("J116").Value = 1095.6666
$ws.Range("K116").Value = 2058
$ws.Range("L116").Value = 1095.6666
$ws.Range("M116").Value = 236
$ws.Range("N116").Value = -5683.6666
$ws.Range("H122").Value = 325835.47
$ws.Range("I122").Value = 502745.16
$ws.Range("K122").Value = 1508235.48
$ws.Range("M122").Value = -1505785.48
$ws.Range("H132").Value = 6609.643
$ws.Range("I132").Value = 2088.1667
$ws.Range("J132").Value = 10000.75
$ws.Range("K132").Value = 6264.500100000001
$ws.Range("L132").Value = 30002.25
$ws.Range("M132").Value = -3734.500100000001
$ws.Range("N132").Value = -35062.25
$ws.Range("H136").Value = 2399.6667
$ws.Range("I136").Value = 2399
$ws.Range("K136").Value = 7197
$ws.Range("M136").Value = -4647

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1817.4166
$ws.Range("I3").Value = 2058
$ws.Range("J3").Value = 1095.6666
$ws.Range("K3").Value = 2058
$ws.Range("L3").Value = 1095.6666
$ws.Range("M3").Value = -1944
$ws.Range("N3").Value = -1323.6666
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 150
$ws.Range("K4").Value = 150
$ws.Range("M4").Value = -35
$ws.Range("H22").Value = 528.3333
$ws.Range("I22").Value = 567.2727
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 567.2727
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = -394.2727
$ws.Range("N22").Value = -446
$ws.Range("H99").Value = 1858.25
$ws.Range("I99").Value = 1754.909
$ws.Range("K99").Value = 1754.909
$ws.Range("M99").Value = -256.9090000000001
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 346.2
$ws.Range("I22").Value = 307.75
$ws.Range("K22").Value = 307.75
$ws.Range("M22").Value = 42.25
$ws.Range("H31").Value = 4981.1333
$ws.Range("I31").Value = 4739.9375
$ws.Range("J31").Value = 5256.7856
$ws.Range("K31").Value = 4739.9375
$ws.Range("L31").Value = 5256.7856
$ws.Range("M31").Value = -4444.9375
$ws.Range("N31").Value = -5846.7856
$ws.Range("H34").Value = 4981.1333
$ws.Range("I34").Value = 4739.9375
$ws.Range("J34").Value = 5256.7856
$ws.Range("K34").Value = 4739.9375
$ws.Range("L34").Value = 5256.7856
$ws.Range("M34").Value = -4537.9375
$ws.Range("N34").Value = -5660.7856
$ws.Range("H58").Value = 4667.619
$ws.Range("I58").Value = 2538.923
$ws.Range("K58").Value = 2538.923
$ws.Range("M58").Value = -2335.923
$ws.Range("H68").Value = 47666.332
$ws.Range("J68").Value = 47666.332
$ws.Range("L68").Value = 47666.332
$ws.Range("N68").Value = -49164.332
$ws.Range("H71").Value = 47666.332
$ws.Range("J71").Value = 47666.332
$ws.Range("L71").Value = 142998.996
$ws.Range("N71").Value = -150486.996
$ws.Range("H86").Value = 12846.75
$ws.Range("I86").Value = 4787.5
$ws.Range("J86").Value = 16876.375
$ws.Range("K86").Value = 4787.5
$ws.Range("L86").Value = 16876.375
$ws.Range("M86").Value = -3664.5
$ws.Range("N86").Value = -19122.375
$ws.Range("H89").Value = 12846.75
$ws.Range("I89").Value = 4787.5
$ws.Range("J89").Value = 16876.375
$ws.Range("K89").Value = 23937.5
$ws.Range("L89").Value = 84381.875
$ws.Range("M89").Value = -18321.5
$ws.Range("N89").Value = -95613.875
$ws.Range("H94").Value = 1403
$ws.Range("J94").Value = 1403
$ws.Range("L94").Value = 1403
$ws.Range("N94").Value = -2305
$ws.Range("H105").Value = 2932.8333
$ws.Range("I105").Value = 2720.125
$ws.Range("J105").Value = 3358.25
$ws.Range("K105").Value = 2720.125
$ws.Range("L105").Value = 3358.25
$ws.Range("M105").Value = -973.125
$ws.Range("N105").Value = -6852.25
$ws.Range("H107").Value = 1079.875
$ws.Range("I107").Value = 546.3333
$ws.Range("K107").Value = 546.3333
$ws.Range("M107").Value = 1373.6667
$ws.Range("H122").Value = 968.25
$ws.Range("I122").Value = 991
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 2973
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -523
$ws.Range("N122").Value = -7600
$ws.Range("H125").Value = 89999.5
$ws.Range("J125").Value = 89999.5
$ws.Range("L125").Value = 89999.5
$ws.Range("N125").Value = -94919.5
$ws.Range("H134").Value = 2688.1365
$ws.Range("I134").Value = 2021.2106
$ws.Range("K134").Value = 6063.6318
$ws.Range("M134").Value = -3528.6318
$ws.Range("H136").Value = 4667.619
$ws.Range("I136").Value = 2538.923
$ws.Range("K136").Value = 7616.768999999999
$ws.Range("M136").Value = -5066.768999999999
$ws.Range("H141").Value = 120714.14
$ws.Range("J141").Value = 120714.14
$ws.Range("L141").Value = 120714.14
$ws.Range("N141").Value = -131074.14

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 3545
$ws.Range("I14").Value = 3545
$ws.Range("K14").Value = 10635
$ws.Range("M14").Value = -10462
$ws.Range("H68").Value = 2907.3635
$ws.Range("J68").Value = 2753.6667
$ws.Range("L68").Value = 8261.000100000001
$ws.Range("N68").Value = -9883.000100000001
$ws.Range("H71").Value = 2907.3635
$ws.Range("J71").Value = 2753.6667
$ws.Range("L71").Value = 24783.0003
$ws.Range("N71").Value = -32895.0003
$ws.Range("H80").Value = 5997.8335
$ws.Range("I80").Value = 5995.6665
$ws.Range("K80").Value = 17986.9995
$ws.Range("M80").Value = -17050.9995
$ws.Range("H81").Value = 2781.3333
$ws.Range("J81").Value = 2781.3333
$ws.Range("L81").Value = 8343.999899999999
$ws.Range("N81").Value = -10589.9999
$ws.Range("H83").Value = 5997.8335
$ws.Range("I83").Value = 5995.6665
$ws.Range("K83").Value = 53960.9985
$ws.Range("M83").Value = -49280.9985
$ws.Range("H84").Value = 2781.3333
$ws.Range("J84").Value = 2781.3333
$ws.Range("L84").Value = 25031.9997
$ws.Range("N84").Value = -36263.9997
$ws.Range("H104").Value = 84327.336
$ws.Range("I104").Value = 1488
$ws.Range("J104").Value = 167166.67
$ws.Range("K104").Value = 4464
$ws.Range("L104").Value = 501500.01
$ws.Range("M104").Value = -1843
$ws.Range("N104").Value = -506742.01
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H113").Value = 3209.6667
$ws.Range("J113").Value = 3678
$ws.Range("L113").Value = 11034
$ws.Range("N113").Value = -15374
$ws.Range("H132").Value = 1998
$ws.Range("J132").Value = 1998
$ws.Range("L132").Value = 17982
$ws.Range("N132").Value = -23042
$ws.Range("N106").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 4732.5
$ws.Range("J55").Value = 5966.6665
$ws.Range("L55").Value = 5966.6665
$ws.Range("N55").Value = -6620.6665
$ws.Range("H102").Value = 3049.3572
$ws.Range("I102").Value = 1965.6666
$ws.Range("K102").Value = 1965.6666
$ws.Range("M102").Value = -343.6666
$ws.Range("H113").Value = 5101.1
$ws.Range("H122").Value = 482107.47
$ws.Range("I122").Value = 62263.707
$ws.Range("K122").Value = 186791.121
$ws.Range("M122").Value = -184341.121
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H126").Value = 4996.5
$ws.Range("I126").Value = 4982.5
$ws.Range("K126").Value = 14947.5
$ws.Range("M126").Value = -12477.5
$ws.Range("H132").Value = 3425.923
$ws.Range("I132").Value = 2444.7
$ws.Range("J132").Value = 6696.6665
$ws.Range("K132").Value = 7334.099999999999
$ws.Range("L132").Value = 20089.9995
$ws.Range("M132").Value = -4804.099999999999
$ws.Range("N132").Value = -25149.9995
$ws.Range("H141").Value = 70824.25
$ws.Range("J141").Value = 70824.25
$ws.Range("L141").Value = 70824.25
$ws.Range("N141").Value = -81184.25
$ws.Range("N124").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2749.75
$ws.Range("J7").Value = 3499.5
$ws.Range("L7").Value = 3499.5
$ws.Range("N7").Value = -3723.5
$ws.Range("H16").Value = 1576.7059
$ws.Range("I16").Value = 1654
$ws.Range("J16").Value = 1216
$ws.Range("K16").Value = 1654
$ws.Range("L16").Value = 1216
$ws.Range("M16").Value = -1484
$ws.Range("N16").Value = -1556
$ws.Range("H22").Value = 2299
$ws.Range("J22").Value = 3274
$ws.Range("L22").Value = 3274
$ws.Range("N22").Value = -3864
$ws.Range("H27").Value = 2299
$ws.Range("J27").Value = 3274
$ws.Range("L27").Value = 3274
$ws.Range("N27").Value = -3488
$ws.Range("H46").Value = 3639.9333
$ws.Range("J46").Value = 4971.2856
$ws.Range("L46").Value = 4971.2856
$ws.Range("N46").Value = -5347.2856
$ws.Range("H55").Value = 877.4231
$ws.Range("I55").Value = 645.2105
$ws.Range("K55").Value = 645.2105
$ws.Range("M55").Value = -472.2105
$ws.Range("H56").Value = 1650
$ws.Range("I56").Value = 1650
$ws.Range("K56").Value = 1650
$ws.Range("M56").Value = -959
$ws.Range("H61").Value = 3677.8147
$ws.Range("I61").Value = 3387.625
$ws.Range("K61").Value = 3387.625
$ws.Range("M61").Value = -3185.625
$ws.Range("H68").Value = 4198.8
$ws.Range("I68").Value = 3998
$ws.Range("K68").Value = 3998
$ws.Range("M68").Value = -3249
$ws.Range("H71").Value = 4198.8
$ws.Range("I71").Value = 3998
$ws.Range("K71").Value = 19990
$ws.Range("M71").Value = -16246
$ws.Range("H93").Value = 1046.7
$ws.Range("I93").Value = 495.2857
$ws.Range("J93").Value = 2333.3333
$ws.Range("K93").Value = 495.2857
$ws.Range("L93").Value = 2333.3333
$ws.Range("M93").Value = 752.7143
$ws.Range("N93").Value = -4829.3333
$ws.Range("H113").Value = 3677.8147
$ws.Range("I113").Value = 3387.625
$ws.Range("K113").Value = 3387.625
$ws.Range("M113").Value = -1217.625
$ws.Range("H122").Value = 7375.9287
$ws.Range("I122").Value = 6256.143
$ws.Range("J122").Value = 8495.714
$ws.Range("K122").Value = 18768.429
$ws.Range("L122").Value = 25487.142
$ws.Range("M122").Value = -16318.429
$ws.Range("N122").Value = -30387.142
$ws.Range("H126").Value = 2749.75
$ws.Range("J126").Value = 3499.5
$ws.Range("L126").Value = 10498.5
$ws.Range("N126").Value = -15438.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 45007.5
$ws.Range("J51").Value = 50015
$ws.Range("L51").Value = 50015
$ws.Range("N51").Value = -51035
$ws.Range("H107").Value = 3333.6667
$ws.Range("I107").Value = 2002
$ws.Range("J107").Value = 3999.5
$ws.Range("K107").Value = 6006
$ws.Range("L107").Value = 11998.5
$ws.Range("M107").Value = -4086
$ws.Range("N107").Value = -15838.5
$ws.Range("H113").Value = 2621
$ws.Range("J113").Value = 3951.5
$ws.Range("L113").Value = 11854.5
$ws.Range("N113").Value = -16194.5
$ws.Range("H122").Value = 2250.5
$ws.Range("I122").Value = 2250.5
$ws.Range("K122").Value = 6751.5
$ws.Range("M122").Value = -4301.5
$ws.Range("H126").Value = 256247.25
$ws.Range("J126").Value = 15000
$ws.Range("L126").Value = 45000
$ws.Range("N126").Value = -49940
